# Update volume tables with revised statistics and values
# (Table2_intraday_volume_total.xlsx) - rows 2 and 3 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - "Ann Window Volume" for Panel A (Bond Futures) / FF1
$ws.Range("D2").Value = 60.40670740372318
$ws.Range("E2").Value = 212.6735049918422
$ws.Range("G2").Value = 0.5483870967741935
$ws.Range("H2").Value = 32.29032258064516
$ws.Range("J2").Value = 54.57889795623601
$ws.Range("K2").Value = 164.6229224345818
$ws.Range("M2").Value = 1.655737704918033
$ws.Range("N2").Value = 28.30327868852459
$ws.Range("P2").Value = 40.75869224888047
$ws.Range("Q2").Value = 109.5695240186682
$ws.Range("S2").Value = 1.628099173553719
$ws.Range("T2").Value = 25.93388429752066
$ws.Range("V2").Value = 28.8698762324313
$ws.Range("W2").Value = 73.36182663070548
$ws.Range("Y2").Value = 2.028571428571428
$ws.Range("Z2").Value = 19.23571428571429
$ws.Range("AB2").Value = 7.631100654118275
$ws.Range("AC2").Value = 17.61769816597677
$ws.Range("AE2").Value = 0.8075757575757576
$ws.Range("AF2").Value = 6.338636363636363

# Row 3 - "Diff (Ann - Non)" for Panel A (Bond Futures) / FF1
$ws.Range("D3").Value = 55.07126616455876
$ws.Range("J3").Value = 49.58332129703184
$ws.Range("P3").Value = 36.29712018058033
$ws.Range("V3").Value = 24.1410950283197
$ws.Range("AB3").Value = 4.954557969563476
